# Hypnos V3 hardware fix - bill of materials update
# Quantity for the part in row 7 (Digikey/Arrow/Mouser sheets) drops from
# 4 to 2; the dependent unit-cost (column F, row 7) and the sheet total
# (column F, row 20) are formulas and recalc automatically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Digikey", "Arrow", "Mouser")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()

    # Quantity changed from 4 to 2
    $ws.Range("C7").Value = 2

    # Cursor ends up on C8 after the edit (matches the saved selection)
    $ws.Range("C8").Select()
}

# Restore the originally active sheet (Mouser, tabSelected in the source file)
$wb.Worksheets.Item("Mouser").Activate()
